$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("102")
$ws2 = $wb.Worksheets.Item("103")

# --- Update data on sheet "103" (rows 8-10 updated, rows 11-12 added) ---

# Row 8: DCLink controller info updated
$ws2.Range("B8").Value = "IA-03RaCtrl:CO-DCLinkCtrl"
$ws2.Range("E8").Value = "IA-03RaPS01:PS-DCLink-SI, IA-03RaPS02:PS-DCLink-SI"

# Row 9: now describes PSCtrl-SI1 (was BBB-SI-CORRETORAS1)
$ws2.Range("A9").Value = "10.128.103.121"
$ws2.Range("B9").Value = "IA-03RaCtrl:CO-PSCtrl-SI1"
$ws2.Range("D9").Value = "1,2,3,4,5,6,7,8"
$ws2.Range("E9").Value = "SI-03M2:PS-QFP, SI-03M2:PS-QDP1, SI-03M2:PS-QDP2, SI-04M1:PS-QFB, SI-04M1:PS-QDB1, SI-04M1:PS-QDB2, SI-03M1:PS-QS, SI-03M2:PS-QS"

# Row 10: now describes PSCtrl-SI2 (was BBB-SI-CORRETORAS2)
$ws2.Range("A10").Value = "10.128.103.122"
$ws2.Range("B10").Value = "IA-03RaCtrl:CO-PSCtrl-SI2"
$ws2.Range("E10").Value = "SI-03M1:PS-CH, SI-03M1:PS-CV, SI-03M2:PS-CH, SI-03M2:PS-CV, SI-03C2:PS-CH, SI-03C2:PS-CV-1, SI-03C2:PS-CV-2"

# Row 11 (new): PSCtrl-SI3
$ws2.Range("A11").Value = "10.128.103.131"
$ws2.Range("B11").Value = "IA-03RaCtrl:CO-PSCtrl-SI3"
$ws2.Range("C11").Value = "PowerSupply"
$ws2.Range("D11").Value = "1,2,3,4,5,6,7,8,9,10,11"
$ws2.Range("E11").Value = "SI-03C1:PS-Q1, SI-03C1:PS-Q2, SI-03C2:PS-Q3, SI-03C2:PS-Q4, SI-03C4:PS-Q1, SI-03C4:PS-Q2, SI-03C3:PS-Q3, SI-03C3:PS-Q4, SI-03C1:PS-QS, SI-03C2:PS-QS, SI-03C3:PS-QS"

# Row 12 (new): PSCtrl-SI4
$ws2.Range("A12").Value = "10.128.103.132"
$ws2.Range("B12").Value = "IA-03RaCtrl:CO-PSCtrl-SI4"
$ws2.Range("C12").Value = "PowerSupply"
$ws2.Range("D12").Value = "1,2,3,4,5,6,7"
$ws2.Range("E12").Value = "SI-03C1:PS-CH, SI-03C1:PS-CV, SI-03C4:PS-CH, SI-03C4:PS-CV, SI-03C3:PS-CH, SI-03C3:PS-CV-1, SI-03C3:PS-CV-2"

# --- View state: sheet "102" loses selection/topLeftCell/tabSelected, "103" becomes active ---
[void]$ws1.Activate()
[void]$ws1.Range("B8").Select()

[void]$ws2.Activate()
[void]$ws2.Range("E13").Select()
